$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 204, pushing existing rows 204-265 down to 206-267
$ws.Rows.Item(204).Insert()
$ws.Rows.Item(204).Insert()

# New row 204
$ws.Range("A204").Value = 4
$ws.Range("B204").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C204").Value = "Los Lagos"
$ws.Range("D204").Value = 44463
$ws.Range("E204").Value = 10
$ws.Range("F204").Value = "Fruta"
$ws.Range("G204").Value = 100102
$ws.Range("H204").Value = "Cítricos"
$ws.Range("I204").Value = 100102003
$ws.Range("J204").Value = "Limón"
$ws.Range("K204").Value = "Sin especificar"
$ws.Range("L204").Value = "1a amarillo"
$ws.Range("M204").Value = 1000
$ws.Range("N204").Value = 10000
$ws.Range("O204").Value = 10000
$ws.Range("P204").Value = 10000
$ws.Range("Q204").Value = "`$/malla 18 kilos"
$ws.Range("R204").Value = "Provincia de Melipilla"
$ws.Range("S204").Value = 556
$ws.Range("T204").Value = 18

# New row 205
$ws.Range("A205").Value = 4
$ws.Range("B205").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C205").Value = "Los Lagos"
$ws.Range("D205").Value = 44463
$ws.Range("E205").Value = 10
$ws.Range("F205").Value = "Fruta"
$ws.Range("G205").Value = 100102
$ws.Range("H205").Value = "Cítricos"
$ws.Range("I205").Value = 100102003
$ws.Range("J205").Value = "Limón"
$ws.Range("K205").Value = "Sin especificar"
$ws.Range("L205").Value = "2a amarillo"
$ws.Range("M205").Value = 500
$ws.Range("N205").Value = 9000
$ws.Range("O205").Value = 9000
$ws.Range("P205").Value = 9000
$ws.Range("Q205").Value = "`$/malla 18 kilos"
$ws.Range("R205").Value = "Provincia de Melipilla"
$ws.Range("S205").Value = 500
$ws.Range("T205").Value = 18
